$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring over the same date/time formatting used by the prior week's row
# before filling in the new values, so the new cells match existing style.
$ws.Range("A8:C8").Copy() | Out-Null
$ws.Range("A9:C9").PasteSpecial(-4122) | Out-Null

# New row of data for the week starting 2014-07-07 (serial 41827)
$ws.Range("A9").Value = 41827
$ws.Range("B9").Value = 0.77083333333333337
$ws.Range("C9").Value = 0.875

# New entry describing the work done, placed in column E of row 8
$ws.Range("E8").Value = "More formatting for the PDF generation, and now its more modular"

# Update the selected cell to reflect where the user left off editing
$ws.Range("E10").Select()

$wb.Save()
